# Update the "ContactUs" test-data sheet with new ClickAction page test case
# data: fill in rows 2-9 (firstname/lastname/email/comment), add a
# Submit Expected/Actual Result column, style + hyperlink the email
# column, and select F7 to match the author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ContactUs")

# ---- Phase 1: fill rows 2-9 for columns A (firstname) - D (comment) ----
$ws.Range("A2").Value = "Bao"
$ws.Range("B2").Value = "Nguyen"
$ws.Range("C2").Value = "abc@gmail.com"
$ws.Range("D2").Value = "Note1"

$ws.Range("B3").Value = "Nguyen"
$ws.Range("C3").Value = "abc@gmail.com"
$ws.Range("D3").Value = "Note2"

$ws.Range("A4").Value = "Bao"
$ws.Range("C4").Value = "abc@gmail.com"
$ws.Range("D4").Value = "Note3"

$ws.Range("C5").Value = "abc@gmail.com"
$ws.Range("D5").Value = "Note4"

$ws.Range("A6").Value = "Bao"
$ws.Range("B6").Value = "Nguyen"
$ws.Range("D6").Value = "Note5"

$ws.Range("A7").Value = "Bao"
$ws.Range("B7").Value = "Nguyen"
$ws.Range("C7").Value = "abc@gmail.com"

$ws.Range("A8").Value = "123&**"
$ws.Range("B8").Value = "Nguyen"
$ws.Range("C8").Value = "abc"
$ws.Range("D8").Value = "Note7"

$ws.Range("A9").Value = 12239
$ws.Range("B9").Value = "*&^^()!"

# ---- Phase 2: fill the new "Submit - Expected Result" column (E2:E9) ----
$ws.Range("E2").Value = "Passed"
$ws.Range("E3").Value = "Failed"
$ws.Range("E4").Value = "Failed"
$ws.Range("E5").Value = "Failed"
$ws.Range("E6").Value = "Failed"
$ws.Range("E7").Value = "Failed"
$ws.Range("E8").Value = "Failed"
$ws.Range("E9").Value = "Failed"

# ---- Phase 3: add the new header labels last ----
$ws.Range("E1").Value = "Submit - Expected Result"
$ws.Range("F1").Value = "Submit - Actual Result"

# ---- Hyperlink the email cells (adds Hyperlink style too) ----
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:abc@gmail.com", "", "", "abc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:abc@gmail.com", "", "", "abc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:abc@gmail.com", "", "", "abc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:abc@gmail.com", "", "", "abc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:abc@gmail.com", "", "", "abc@gmail.com")

# ---- Column widths (best-fit sized, like an AutoFit after data entry) ----
$ws.Columns.Item(1).ColumnWidth = 7.944010416666667
$ws.Columns.Item(2).ColumnWidth = 7.608072916666667
$ws.Columns.Item(3).ColumnWidth = 13.385416666666666
$ws.Columns.Item(4).ColumnWidth = 7.944010416666667
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668
$ws.Columns.Item(6).ColumnWidth = 21.498697916666668

# ---- Final selection state ----
$ws.Range("F7").Select()
